$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cell A1
$ws.Range("A1").Value = "input_keyName"

# Remove column B entirely (its header/value cells go away and the grid
# collapses back down to a single column)
$ws.Columns.Item(2).Delete()

# Widen column A to 15 characters
$ws.Columns.Item(1).ColumnWidth = 14.1
